# Add an "ara" (Arabic) row to the status_type master-data sheet, mirroring
# the existing "fra" row's layout but with Arabic text wrapped in C/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- New row values -----------------------------------------------------
$ws.Cells.Item(4, 1).Value = "ara"
$ws.Cells.Item(4, 2).Value = "ADT"

$arabicText = "حالة التفعيل أو التشغيل"
$ws.Cells.Item(4, 3).Value = $arabicText
$ws.Cells.Item(4, 4).Value = $arabicText

# E4 should hold the literal text "TRUE" (shared string), same as E2/E3 -
# copy it from E2 so it lands as text rather than being auto-coerced to a
# boolean.
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(4, 5).PasteSpecial(-4163)
$excel.CutCopyMode = $false

# -- Formatting for the new Arabic cells (wrap + left align) ------------
$c4 = $ws.Cells.Item(4, 3)
$c4.WrapText = $true
$c4.HorizontalAlignment = -4131

$c4.Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -- Row height / column widths ------------------------------------------
$ws.Rows.Item(4).RowHeight = 16.4
$ws.Columns.Item(3).ColumnWidth = 22.6
$ws.Columns.Item(4).ColumnWidth = 19.95

# -- Selection, matching the authored diff -------------------------------
$ws.Range("C4:D4").Select()
